$d = $word.ActiveDocument

# Start from the end of the last paragraph in the document body (before sectPr)
$lastPara = $d.Paragraphs.Last
$rng = $lastPara.Range
$rng.Collapse(0)

# --- New bullet (ilvl 0): "Pobreza_Monetaria" ---
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$p1 = $d.Paragraphs.Last
$p1.Range.ListFormat.ListLevelNumber = 1
$p1.Range.Text = "Pobreza_Monetaria"

# --- New sub-bullet (ilvl 1): description ---
$rng2 = $p1.Range
$rng2.Collapse(0)
$rng2.InsertParagraphAfter()
$rng2.Collapse(0)
$p2 = $d.Paragraphs.Last
$p2.Range.ListFormat.ListLevelNumber = 2
$p2.Range.Text = "The information from this dataset is stored at the state level, we are going to assume the values analyzed here have the same values as their largest city (aka the capital city of each of the states). "
